$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" cells get a value that Excel would otherwise auto-detect
# as a number (e.g. "0.999", "219.49"). The source workbook stores every
# Price/Volume cell as text, so force a text format on just those cells
# before assigning the value to keep them as text, matching the workbook.
$textCells = @("D4", "D5", "D8", "D10", "D11", "D14", "D15", "D16", "D18", "D19", "D21", "D24", "D25", "D26", "D28", "D29", "D30", "D35", "D36", "D38", "D43", "D44", "D45", "D46", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.859.11"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "1.678.18"
$ws.Range("E3").Value = "  +3.12%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "219.49"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "28.99"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "0.0642"
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.921.80"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "1.682.53"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.602"
$ws.Range("E14").Value = "  +7.18%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "10.07"
$ws.Range("E15").Value = "  +9.33%  "
$ws.Range("D16").Value = "4.05"
$ws.Range("E16").Value = "  +5.80%  "
$ws.Range("D17").Value = "30.845.35"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "65.89"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").Value = "242.93"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "0.0₃0718"
$ws.Range("E20").Value = "  +2.51%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "158.94"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "15.81"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").Value = "6.67"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "0.0492"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").Value = "1.514.60"
$ws.Range("E33").Value = "  +6.52%  "
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("D35").Value = "1.74"
$ws.Range("E35").Value = "  +4.37%  "
$ws.Range("D36").Value = "83.92"
$ws.Range("E36").Value = "  +12.77%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "0.607"
$ws.Range("E38").Value = "  +9.31%  "
$ws.Range("E39").Value = "  +5.11%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("E42").Value = "  +3.18%  "
$ws.Range("D43").Value = "0.838"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "0.0500"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "1.03"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +4.53%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.813.57"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "50.66"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("E50").Value = "  +6.39%  "
$ws.Range("D51").Value = "92.83"
$ws.Range("E51").Value = "  +2.14%  "
